$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -7.995
    18 = -8.377000000000001
    20 = -7.691000000000001
    27 = -7.866000000000001
    35 = -7.888
    69 = -7.679
    76 = -7.943000000000001
    78 = -8.077
    82 = -8.170000000000002
    83 = -8.103
    93 = -6.976999999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
